# Update countries & provincias Spain
#
# 1) Refresh COVID stats for several existing countries (Estados Unidos,
#    Austria, Luxemburgo, Argelia, Eslovaquia, San Marino).
# 2) Insert a new "Congo" entry (with its own stats) in its sorted position
#    (between Guadalupe and Jamaica) and remove the old "Congo" entry that
#    used to sit between Camboya and Madagascar.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Simple stat refreshes -------------------------------------------

# Row 4: Estados Unidos
$ws.Cells.Item(4,2).Value = 684427
$ws.Cells.Item(4,3).Value = 6857
$ws.Cells.Item(4,4).Value = 58156
$ws.Cells.Item(4,5).Value = 590808
$ws.Cells.Item(4,7).Value = 846
$ws.Cells.Item(4,8).Value = 35463

# Row 20: Austria
$ws.Cells.Item(20,2).Value = 14570
$ws.Cells.Item(20,3).Value = 94
$ws.Cells.Item(20,5).Value = 4456

# Row 51: Luxemburgo
$ws.Cells.Item(51,2).Value = 3480
$ws.Cells.Item(51,3).Value = 36
$ws.Cells.Item(51,4).Value = 579
$ws.Cells.Item(51,5).Value = 2829
$ws.Cells.Item(51,6).Value = 29
$ws.Cells.Item(51,7).Value = 3
$ws.Cells.Item(51,8).Value = 72

# Row 58: Argelia
$ws.Cells.Item(58,2).Value = 2418
$ws.Cells.Item(58,3).Value = 150
$ws.Cells.Item(58,4).Value = 846
$ws.Cells.Item(58,5).Value = 1208
$ws.Cells.Item(58,7).Value = 16
$ws.Cells.Item(58,8).Value = 364

# Row 79: Eslovaquia
$ws.Cells.Item(79,4).Value = 175
$ws.Cells.Item(79,5).Value = 865
$ws.Cells.Item(79,7).Value = 1
$ws.Cells.Item(79,8).Value = 9

# Row 104: San Marino
$ws.Cells.Item(104,2).Value = 435
$ws.Cells.Item(104,3).Value = 9
$ws.Cells.Item(104,4).Value = 57
$ws.Cells.Item(104,5).Value = 339
$ws.Cells.Item(104,7).Value = 1
$ws.Cells.Item(104,8).Value = 39

# --- 2) Move "Congo" to its new sorted position -------------------------

# Insert a new row right above "Jamaica" (row 129) to hold Congo's new
# data. This shifts Jamaica, Ruanda, Brunei, Gibraltar and Camboya down by
# one row each, keeping their own data intact.
$ws.Rows.Item(129).Insert()

$ws.Cells.Item(129,1).Value = "Congo"
$ws.Cells.Item(129,2).Value = 143
$ws.Cells.Item(129,3).Value = 26
$ws.Cells.Item(129,4).Value = 11
$ws.Cells.Item(129,5).Value = 126
$ws.Cells.Item(129,6).Value = 0
$ws.Cells.Item(129,7).Value = 1
$ws.Cells.Item(129,8).Value = 6

# The old "Congo" row (previously row 134, now pushed down to row 135 by
# the insert above) is now obsolete; remove it so Madagascar and the rest
# of the list shift back up into their original rows.
$ws.Rows.Item(135).Delete()
